$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price strings that look numeric (e.g. "1.002") stay as plain text,
# matching how the source data was scraped/stored (inline strings).
$ws.Range("D4:D5").NumberFormat = "@"
$ws.Range("D7:D16").NumberFormat = "@"
$ws.Range("D18:D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25:D31").NumberFormat = "@"
$ws.Range("D33:D47").NumberFormat = "@"
$ws.Range("D49:D51").NumberFormat = "@"

# Update coin rows with refreshed price / volume data
$ws.Range("D2").Value = "24.569.32"
$ws.Range("E2").Value = "  +3.46%  "

$ws.Range("D3").Value = "1.694.11"
$ws.Range("E3").Value = "  +1.73%  "

$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").Value = "315.91"
$ws.Range("E5").Value = "  +2.07%  "

$ws.Range("E6").Value = "  +0.12%  "

$ws.Range("D7").Value = "0.3936"
$ws.Range("E7").Value = "  +1.23%  "

$ws.Range("D8").Value = "0.4010"
$ws.Range("E8").Value = "  +1.75%  "

$ws.Range("D9").Value = "1.524"
$ws.Range("E9").Value = "  +5.99%  "

$ws.Range("D10").Value = "1.002"
$ws.Range("E10").Value = "  +0.22%  "

$ws.Range("D11").Value = "52.80"
$ws.Range("E11").Value = "  +6.15%  "

$ws.Range("D12").Value = "0.08736"
$ws.Range("E12").Value = "  +0.90%  "

$ws.Range("D13").Value = "7.199"
$ws.Range("E13").Value = "  +7.13%  "

$ws.Range("D14").Value = "23.12"
$ws.Range("E14").Value = "  +2.32%  "

$ws.Range("D15").Value = "0.00001317"
$ws.Range("E15").Value = "  +0.47%  "

$ws.Range("D16").Value = "7.559"
$ws.Range("E16").Value = "  +4.22%  "

$ws.Range("D17").Value = "1.693.78"
$ws.Range("E17").Value = "  +1.62%  "

$ws.Range("D18").Value = "99.72"
$ws.Range("E18").Value = "  +0.21%  "

$ws.Range("D19").Value = "0.07049"
$ws.Range("E19").Value = "  +3.85%  "

$ws.Range("E20").Value = "  +3.36%  "

$ws.Range("D21").Value = "6.864"
$ws.Range("E21").Value = "  +3.54%  "

$ws.Range("E22").Value = "  +0.03%  "

$ws.Range("D23").Value = "14.03"
$ws.Range("E23").Value = "  +1.62%  "

$ws.Range("D24").Value = "24.554.42"
$ws.Range("E24").Value = "  +3.41%  "

$ws.Range("D25").Value = "3.005"
$ws.Range("E25").Value = "  +6.47%  "

$ws.Range("D26").Value = "2.321"
$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("D27").Value = "22.34"
$ws.Range("E27").Value = "  +2.90%  "

$ws.Range("D28").Value = "160.20"
$ws.Range("E28").Value = "  +0.69%  "

$ws.Range("D29").Value = "5.219"
$ws.Range("E29").Value = "  +1.31%  "

$ws.Range("D30").Value = "134.49"
$ws.Range("E30").Value = "  +3.83%  "

$ws.Range("D31").Value = "7.423"
$ws.Range("E31").Value = "  +9.90%  "

$ws.Range("D32").Value = "1.880.88"
$ws.Range("E32").Value = "  +1.73%  "

$ws.Range("D33").Value = "1.096"
$ws.Range("E33").Value = "  -2.13%  "

$ws.Range("D34").Value = "0.08507"
$ws.Range("E34").Value = "  -0.01%  "

$ws.Range("D35").Value = "7.209"
$ws.Range("E35").Value = "  +7.62%  "

$ws.Range("D36").Value = "11.45"
$ws.Range("E36").Value = "  +9.51%  "

$ws.Range("D37").Value = "1.957"
$ws.Range("E37").Value = "  +0.21%  "

$ws.Range("D38").Value = "0.2723"
$ws.Range("E38").Value = "  +2.37%  "

$ws.Range("D39").Value = "14.42"
$ws.Range("E39").Value = "  -0.04%  "

$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").Value = "0.09046"
$ws.Range("E40").Value = "  +2.94%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "0.02737"
$ws.Range("E41").Value = "  +9.10%  "

$ws.Range("D42").Value = "1.461"
$ws.Range("E42").Value = "  +0.68%  "

$ws.Range("D43").Value = "0.7682"
$ws.Range("E43").Value = "  +2.17%  "

$ws.Range("D44").Value = "0.7190"
$ws.Range("E44").Value = "  +2.73%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "15.48"
$ws.Range("E45").Value = "  +4.65%  "

$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "2.537"
$ws.Range("E46").Value = "  +5.71%  "

$ws.Range("D47").Value = "4.207"
$ws.Range("E47").Value = "  +2.54%  "

$ws.Range("E48").Value = "  +0.10%  "

$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "140.85"
$ws.Range("E49").Value = "  +1.53%  "

$ws.Range("B50").Value = "Flow"
$ws.Range("C50").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D50").Value = "1.322"
$ws.Range("E50").Value = "  +9.82%  "

$ws.Range("D51").Value = "0.08001"
$ws.Range("E51").Value = "  +3.14%  "

